# Leave Card update — 4/12/2023 4:43 PM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New monthly leave-card rows for 2023-2025 (column PERIOD = first of month) ---
$ws.Range("A81").Value = [DateTime]"2023-01-01"
$ws.Range("C81").Value = 1.25

$ws.Range("A82").Value = [DateTime]"2023-02-01"
$ws.Range("C82").Value = 1.25

$ws.Range("A83").Value = [DateTime]"2023-03-01"
$ws.Range("B83").Value = "VL(3-0-0)"
$ws.Range("C83").Value = 1.25
$ws.Range("D83").Value = 3
$ws.Range("K83").Value = "3/16-17,20/2023"

# --- Header: UNIT now recorded as "CENRO" (moved out of the folder path) ---
$ws.Range("F4").Value = "CENRO"

$ws.Range("A84").Value = [DateTime]"2023-04-01"
$ws.Range("A85").Value = [DateTime]"2023-05-01"
$ws.Range("A86").Value = [DateTime]"2023-06-01"
$ws.Range("A87").Value = [DateTime]"2023-07-01"
$ws.Range("A88").Value = [DateTime]"2023-08-01"
$ws.Range("A89").Value = [DateTime]"2023-09-01"
$ws.Range("A90").Value = [DateTime]"2023-10-01"
$ws.Range("A91").Value = [DateTime]"2023-11-01"
$ws.Range("A92").Value = [DateTime]"2023-12-01"
$ws.Range("A93").Value = [DateTime]"2024-01-01"
$ws.Range("A94").Value = [DateTime]"2024-02-01"
$ws.Range("A95").Value = [DateTime]"2024-03-01"
$ws.Range("A96").Value = [DateTime]"2024-04-01"
$ws.Range("A97").Value = [DateTime]"2024-05-01"
$ws.Range("A98").Value = [DateTime]"2024-06-01"
$ws.Range("A99").Value = [DateTime]"2024-07-01"
$ws.Range("A100").Value = [DateTime]"2024-08-01"
$ws.Range("A101").Value = [DateTime]"2024-09-01"
$ws.Range("A102").Value = [DateTime]"2024-10-01"
$ws.Range("A103").Value = [DateTime]"2024-11-01"
$ws.Range("A104").Value = [DateTime]"2024-12-01"
$ws.Range("A105").Value = [DateTime]"2025-01-01"
$ws.Range("A106").Value = [DateTime]"2025-02-01"
$ws.Range("A107").Value = [DateTime]"2025-03-01"
$ws.Range("A108").Value = [DateTime]"2025-04-01"
$ws.Range("A109").Value = [DateTime]"2025-05-01"
$ws.Range("A110").Value = [DateTime]"2025-06-01"
$ws.Range("A111").Value = [DateTime]"2025-07-01"

# --- Restore view/selection state as closely as COM allows ---
$ws.Activate()
$ws.Range("B3").Select()
$ws.Range("B17").Select()
